# Rename the "Patient" worksheet to "Patients" and update the validation
# message on its header row to match the new tab name. Also make the
# Patients sheet the active tab/selection, as in the authored edit.

$wb = $excel.ActiveWorkbook

$wsPatients = $wb.Worksheets.Item("Patient")
$wsPatients.Name = "Patients"

# Update the header cell that states the required tab name.
$wsPatients.Range("C1").Value = "This tab must be called 'Patients'"

# Make the Patients sheet the active sheet/tab and set its selection.
$wsPatients.Activate()
$wsPatients.Range("B4").Select()
